# New crime data collected - weekly CompStat update (60th Precinct)
# Updates the report title (volume/date range) and the week-to-date /
# 28-day / year-to-date / 2-year crime-complaint figures for rows 14-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Title block: bump the bulletin volume number and the covered-week dates
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 29   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/21/2022  Through  11/27/2022"

# ---------------------------------------------------------------------
# Row 14 (Murder)
# ---------------------------------------------------------------------
$ws.Range("L14").Value = -50

# ---------------------------------------------------------------------
# Row 15 (Rape) - Week-to-date counts are no longer numeric (n/a markers)
# ---------------------------------------------------------------------
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("L15").Value = -35

# ---------------------------------------------------------------------
# Row 16 (Robbery)
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 165
$ws.Range("J16").Value = 98
$ws.Range("K16").Value = 68.367346938775
$ws.Range("L16").Value = 47.321428571428
$ws.Range("M16").Value = -34.262948207171
$ws.Range("N16").Value = -82.972136222910

# ---------------------------------------------------------------------
# Row 17 (Fel. Assault)
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 345
$ws.Range("J17").Value = 319
$ws.Range("K17").Value = 8.150470219435
$ws.Range("L17").Value = 64.285714285714
$ws.Range("M17").Value = 77.835051546391
$ws.Range("N17").Value = -49.413489736070

# ---------------------------------------------------------------------
# Row 18 (Burglary)
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = -30.769230769230
$ws.Range("I18").Value = 147
$ws.Range("J18").Value = 157
$ws.Range("K18").Value = -6.369426751592
$ws.Range("L18").Value = 8.088235294117
$ws.Range("M18").Value = -3.921568627450
$ws.Range("N18").Value = -83.774834437086

# ---------------------------------------------------------------------
# Row 19 (Gr. Larceny)
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 39
$ws.Range("H19").Value = 8.333333333333
$ws.Range("I19").Value = 475
$ws.Range("J19").Value = 372
$ws.Range("K19").Value = 27.688172043010
$ws.Range("L19").Value = 57.284768211920
$ws.Range("M19").Value = -7.945736434108
$ws.Range("N19").Value = -24.841772151898

# ---------------------------------------------------------------------
# Row 20 (G.L.A.)
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 108
$ws.Range("J20").Value = 66
$ws.Range("K20").Value = 63.636363636363
$ws.Range("L20").Value = 74.193548387096
$ws.Range("M20").Value = 5.882352941176
$ws.Range("N20").Value = -89.972144846796

# ---------------------------------------------------------------------
# Row 21 (TOTAL)
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -3.225806451612
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 89
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1258
$ws.Range("J21").Value = 1030
$ws.Range("K21").Value = 22.135922330097
$ws.Range("L21").Value = 47.652582159624
$ws.Range("M21").Value = 1.533494753833
$ws.Range("N21").Value = -70.886368896088

# ---------------------------------------------------------------------
# Row 22 (Transit)
# ---------------------------------------------------------------------
$ws.Range("F22").Value = 3
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 27
$ws.Range("K22").Value = 125
$ws.Range("L22").Value = -10
$ws.Range("M22").Value = -6.896551724137

# ---------------------------------------------------------------------
# Row 23 (Housing)
# ---------------------------------------------------------------------
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 150
$ws.Range("F23").Value = 18
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 163
$ws.Range("J23").Value = 167
$ws.Range("K23").Value = -2.395209580838
$ws.Range("L23").Value = 55.238095238095
$ws.Range("M23").Value = 45.535714285714

# ---------------------------------------------------------------------
# Row 24 (Petit Larceny)
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 20
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 18.918918918918
$ws.Range("I24").Value = 1120
$ws.Range("J24").Value = 795
$ws.Range("K24").Value = 40.880503144654
$ws.Range("L24").Value = 19.148936170212
$ws.Range("M24").Value = 2.096627164995

# ---------------------------------------------------------------------
# Row 25 (Misd. Assault)
# ---------------------------------------------------------------------
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 40
$ws.Range("F25").Value = 34
$ws.Range("H25").Value = -19.047619047619
$ws.Range("I25").Value = 515
$ws.Range("J25").Value = 415
$ws.Range("K25").Value = 24.096385542168
$ws.Range("L25").Value = 30.050505050505
$ws.Range("M25").Value = 6.404958677685

# ---------------------------------------------------------------------
# Row 26 (UCR Rape*) - Week-to-date counts become n/a markers too
# ---------------------------------------------------------------------
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("L26").Value = -7.142857142857

# ---------------------------------------------------------------------
# Row 27 (Other Sex Crimes) - Week-to-date count C27 becomes numeric
# ---------------------------------------------------------------------
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 54
$ws.Range("J27").Value = 33
$ws.Range("K27").Value = 63.636363636363
$ws.Range("L27").Value = 80

# ---------------------------------------------------------------------
# Row 28 (Shooting Vic.)
# ---------------------------------------------------------------------
$ws.Range("L28").Value = 52.941176470588
$ws.Range("M28").Value = 30

# ---------------------------------------------------------------------
# Row 29 (Shooting Inc.)
# ---------------------------------------------------------------------
$ws.Range("L29").Value = 41.666666666666
$ws.Range("M29").Value = 6.25

# ---------------------------------------------------------------------
# Row 30 (Hate Crimes) - 28-day count F30 becomes an n/a marker
# ---------------------------------------------------------------------
$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "0"

# ---------------------------------------------------------------------
# Fix-up pass: restore the original cell styles that NumberFormat/Value
# assignment above disturbed, by pasting formats only from a same-style
# neighbour that was not otherwise touched by this edit.
# ---------------------------------------------------------------------
$ws.Range("A15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

$ws.Range("A26").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("A26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null

$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null

$ws.Range("G30").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
